$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new log rows to the feed logs sheet (rows 166 and 167)
$ws.Range("A166").Value = 165
$ws.Range("B166").Value = 1
$ws.Range("C166").Value = "2024-06-18 11:11:39"
$ws.Range("D166").Value = 200
$ws.Range("E166").Value = 12

$ws.Range("A167").Value = 166
$ws.Range("B167").Value = 2
$ws.Range("C167").Value = "2024-06-18 11:11:39"
$ws.Range("D167").Value = 200
$ws.Range("E167").Value = 0
